$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test data rows (22-30), continuing the existing
# regcntr_id / usr_id / machine_id / lang_code / is_active / cr_by / cr_dtimes pattern.
$newRows = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$r = 22
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $r++
}

# Mirror the saved selection state (user had scrolled/selected down to row 31
# and beyond when the file was last saved).
$ws.Range("A31:XFD1048576").Select()

# Page setup was touched (portrait orientation) before saving.
$ws.PageSetup.Orientation = 1
